# Updating excel output file
# - Rename "Obs" header to "Obs_relatorio"
# - Add a new "Obs_sped" column
# - Replace "VERDADEIRO" validation values with the full validation message
# - Leave the new "Obs_sped" column blank for every data row

$wb = $excel.ActiveWorkbook

$longMsg = "Validado com sucesso! Nenhuma divergência entre o SPED e o relatório foi encontrada!"

# ---- Sheet "Bico": header row 1, data rows 2-15, Obs column H, new column I ----
$ws1 = $wb.Worksheets.Item("Bico")

$ws1.Range("H1").Value = "Obs_relatorio"
$ws1.Range("I1").Value = "Obs_sped"

for ($r = 2; $r -le 15; $r++) {
    $ws1.Cells.Item($r, 8).Value = $longMsg
    $ws1.Cells.Item($r, 9).Value = ""
}

# ---- Sheet "Tanque": header row 1, data rows 2-5, Obs column F, new column G ----
$ws2 = $wb.Worksheets.Item("Tanque")

$ws2.Range("F1").Value = "Obs_relatorio"
$ws2.Range("G1").Value = "Obs_sped"

for ($r = 2; $r -le 5; $r++) {
    $ws2.Cells.Item($r, 6).Value = $longMsg
    $ws2.Cells.Item($r, 7).Value = ""
}
